# Daily attendance processing - 2025-11-26 13:42:26
# Normalizes the "Recorded By" column (G): when the comma-separated list of
# recorders has "System" as either the first or last entry, swap the first
# and last entries so the (properly-cased) "System" entry comes first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $tokens = $text -split ", "
    if ($tokens.Length -lt 2) { continue }

    $first = $tokens[0]
    $last = $tokens[$tokens.Length - 1]

    if (($first.ToLower() -eq "system") -or ($last.ToLower() -eq "system")) {
        $tmp = $tokens[0]
        $tokens[0] = $tokens[$tokens.Length - 1]
        $tokens[$tokens.Length - 1] = $tmp
        $newText = [string]::Join(", ", $tokens)
        $cell.Value2 = $newText
    }
}
